$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 41

$ws.Cells.Item($row, 1).Value = 40
$ws.Cells.Item($row, 2).Value = "india"
$ws.Cells.Item($row, 3).Value = "isl"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45260.64583333334
$ws.Cells.Item($row, 6).Value = "Bengaluru FC"
$ws.Cells.Item($row, 7).Value = 3
$ws.Cells.Item($row, 8).Value = "Punjab"
$ws.Cells.Item($row, 9).Value = 3
$ws.Cells.Item($row, 10).Value = 1.64
$ws.Cells.Item($row, 11).Value = "28/11/2023 15:42"
$ws.Cells.Item($row, 12).Value = 1.76
$ws.Cells.Item($row, 13).Value = "30/11/2023 15:29"
$ws.Cells.Item($row, 14).Value = 3.87
$ws.Cells.Item($row, 15).Value = "28/11/2023 15:42"
$ws.Cells.Item($row, 16).Value = 3.79
$ws.Cells.Item($row, 17).Value = "30/11/2023 15:29"
$ws.Cells.Item($row, 18).Value = 5.24
$ws.Cells.Item($row, 19).Value = "28/11/2023 15:42"
$ws.Cells.Item($row, 20).Value = 4.59
$ws.Cells.Item($row, 21).Value = "30/11/2023 15:29"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/india/isl/bengaluru-fc-minerva-punjab/GdKvj9p9/"

# Copy style from row 40 (previous data row) to keep formatting consistent
$ws.Range("A40:V40").Copy()
$ws.Range("A41:V41").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
